# Add six new article rows (206-211) to the articulos sheet, matching the
# "Add more fields and success message" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 206; Id = 1160; Articulo = "bonito";        Descripcion = "muy bonito"; Coste = 4; Precio = 5 },
    @{ Row = 207; Id = 1161; Articulo = "masbonito";      Descripcion = "muy bonito"; Coste = 4; Precio = 5 },
    @{ Row = 208; Id = 1162; Articulo = "aunmasbonito";   Descripcion = "muy bonito"; Coste = 4; Precio = 5 },
    @{ Row = 209; Id = 1163; Articulo = "superbonito";    Descripcion = "muy bonito"; Coste = 4; Precio = 5 },
    @{ Row = 210; Id = 1164; Articulo = "superbonito2";   Descripcion = "muy bonito"; Coste = 4; Precio = 5 },
    @{ Row = 211; Id = 1165; Articulo = "pepino";         Descripcion = "";           Coste = 0; Precio = 0 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
    $ws.Cells.Item($r.Row, 2).Value = $r.Articulo
    $ws.Cells.Item($r.Row, 3).Value = $r.Descripcion
    $ws.Cells.Item($r.Row, 4).Value = $r.Coste
    $ws.Cells.Item($r.Row, 5).Value = $r.Precio
}

Write-Host "Articulos actualizados correctamente: se agregaron 6 nuevos registros."
